$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New column H: "Impact carbone" header + per-vegetable values.
$ws.Range("H1").Value = "Impact carbone"

$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 2
$ws.Range("H4").Value = 3
$ws.Range("H5").Value = 4
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("H16").Value = 1
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
$ws.Range("H22").Value = 1
$ws.Range("H23").Value = 1
$ws.Range("H24").Value = 1
$ws.Range("H25").Value = 1
$ws.Range("H26").Value = 1
$ws.Range("H27").Value = 1
$ws.Range("H28").Value = 1
$ws.Range("H29").Value = 1
$ws.Range("H30").Value = 1

# Selection moved to J8 on Sheet1 in the saved file.
[void]$ws.Range("J8").Select()
